# EI Variable Installments T2 scenarios
# Insert a new test-step row ("waittopageload1" / 2000) into the
# "Edit Repayment Schedule" sheet just above the existing "clickonsubmit"
# step, and make that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row above the old row 6 ("clickonsubmit"), shifting rows
# 6-12 down to 7-13. The new row inherits the format of the row above it.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6
$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# B6 should carry the same "numeric input" formatting used by the other
# numeric step values (B3/B4), not the text formatting it inherited.
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the selection on this sheet to the newly-added row
$ws.Range("A6:B6").Select()

# Make "Edit Repayment Schedule" the active sheet/tab (was "NewLoanInput")
$ws.Activate()

$wb.Save()
